$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# C1: "line_num_start" -> "idx_start"
$ws.Range("C1").Value = "idx_start"

# E1: new "docstring" header, formatted like the other header cells
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "docstring"

# --- idx_start values shift down by 1 ---
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 44

# --- New docstring column values ---
# Leading apostrophes must stay as literal text (VBA comment syntax), not be
# interpreted as Excel's "quote prefix" cell format, so we build each string
# via a formula (CHAR(39) for the apostrophe) and then convert the formula
# result to a plain value in place.
$ws.Range("E2").Formula = "=CHAR(39)&"" A docstring for a procedure ""&CHAR(10)&CHAR(39)&"" ""&CHAR(10)&CHAR(39)&"" JDL 12/13/21   Modified: 8/1/23 JDL ""&CHAR(10)&CHAR(39)"
$ws.Range("E3").Formula = "=CHAR(39)&"" Method1 docstring is ""&CHAR(10)&CHAR(39)&"" multiline ""&CHAR(10)&CHAR(39)&"" ""&CHAR(10)&CHAR(39)&"" JDL 8/1/23 ""&CHAR(10)&CHAR(39)"
$ws.Range("E4").Formula = "=CHAR(39)&"" Method2 docstring ""&CHAR(10)&CHAR(39)&"" ""&CHAR(10)&CHAR(39)&"" JDL 8/1/23 ""&CHAR(10)&CHAR(39)"

$ws.Range("E2:E4").Copy()
$ws.Range("E2:E4").PasteSpecial(-4163)
